# Atualizando a base de dados com o mes de agosto
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Remove the existing data rows (2-14, newest-first) so that we can
#    rewrite them from scratch in ascending (oldest-first) order together
#    with the new August/2021 draw, without dragging along the old
#    per-row height / per-cell font overrides that used to live on rows 2-3.
# ---------------------------------------------------------------------------
$ws.Rows("2:14").Delete() | Out-Null

# date serial number (1900 date system) -> "dezenas sorteadas" string, in
# chronological order (oldest draw first, newest draw last)
$data = @(
    @{ Date = 44013; Dezenas = "04 10 12 14 36 46" },
    @{ Date = 44044; Dezenas = "09 15 20 33 41 43" },
    @{ Date = 44075; Dezenas = "12 21 29 54 56 57" },
    @{ Date = 44105; Dezenas = "06 07 28 42 45 49" },
    @{ Date = 44136; Dezenas = "02 05 10 29 34 41" },
    @{ Date = 44166; Dezenas = "17 20 22 35 41 42" },
    @{ Date = 44197; Dezenas = "16 21 28 41 49 51" },
    @{ Date = 44228; Dezenas = "02 03 07 48 51 54" },
    @{ Date = 44256; Dezenas = "19 28 30 34 40 51" },
    @{ Date = 44287; Dezenas = "04 27 33 35 38 41" },
    @{ Date = 44317; Dezenas = "12 14 17 18 19 22" },
    @{ Date = 44348; Dezenas = "11 13 16 35 49 50" },
    @{ Date = 44378; Dezenas = "04 11 12 44 45 57" }
)

# newly drawn entry for August/2021, appended after the filter/sort below
$newEntry = @{ Date = 44409; Dezenas = "01 19 35 40 47 54" }

$row = 2
foreach ($item in $data) {
    $a = $ws.Cells.Item($row, 1)
    $a.Value = $item.Date
    $a.NumberFormat = "mmm-yy"

    $b = $ws.Cells.Item($row, 2)
    $b.Value = $item.Dezenas

    $row = $row + 1
}

# ---------------------------------------------------------------------------
# 2. Re-apply the autofilter / sort state on the A1:B14 range *before* the
#    new August draw is appended, so the filter/sort range stays limited to
#    the previously existing rows, exactly like in the source workbook.
# ---------------------------------------------------------------------------
$sortRange = $ws.Range("A1:B14")
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("A1:A14")) | Out-Null
$ws.Sort.SetRange($sortRange)
$ws.Sort.Header = 1
$ws.Sort.Apply()
$sortRange.AutoFilter() | Out-Null

$filterDatabaseName = $ws.Names.Add("_xlnm._FilterDatabase", "=Planilha1!`$A`$1:`$B`$14")
$filterDatabaseName.Visible = $false

# ---------------------------------------------------------------------------
# 2b. Now append the new August/2021 draw as row 15, outside the
#     filter/sort range.
# ---------------------------------------------------------------------------
$newRow = $ws.Cells.Item(15, 1)
$newRow.Value = $newEntry.Date
$newRow.NumberFormat = "mmm-yy"
$ws.Cells.Item(15, 2).Value = $newEntry.Dezenas

# ---------------------------------------------------------------------------
# 3. Refresh the active selection (also drops the stale topLeftCell scroll
#    position left over from before the edit).
# ---------------------------------------------------------------------------
$ws.Range("B9").Select() | Out-Null

# ---------------------------------------------------------------------------
# 4. Misc workbook level tweaks captured by the commit.
# ---------------------------------------------------------------------------
$excel.Calculation = -4135  # xlCalculationManual

$wb.Save()
